# v0.92.5 - Migration Cluster Readiness
# Adds a new "Migration Cluster Readiness" field to the Cluster data
# collection template's header row, alongside Department/Migration
# Cluster/Domain, as part of standardizing the record layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell in column E, row 8 (same row as Department/Department
# Simple/Domain/Migration Cluster)
$ws.Range("E8").Value = "Migration Cluster Readiness"

# Widen column E so the longer header text fits (column previously held an
# unused width left over from an earlier column)
$ws.Columns.Item(5).ColumnWidth = 24.1640625

# Active selection now covers the template's header block, A1:E9
$ws.Range("A1:E9").Select()
